$wb = $excel.ActiveWorkbook

# Rename the existing sheet
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "carsSheet"

# Add a new sheet right after the cars sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "laptopSheet"

# Header row
$ws2.Range("A1").Value = "Brand"
$ws2.Range("B1").Value = "Model"
$ws2.Range("C1").Value = "Color"
$ws2.Range("D1").Value = "Year"
$ws2.Range("A1:D1").Font.Bold = $true

# Data rows
$ws2.Range("A2").Value = "Dell"
$ws2.Range("B2").Value = "Inspiron 15"
$ws2.Range("C2").Value = "Silver"
$ws2.Range("D2").Value = 2022

$ws2.Range("A3").Value = "HP"
$ws2.Range("B3").Value = "Pavilion x360"
$ws2.Range("C3").Value = "Blue"
$ws2.Range("D3").Value = 2021

$ws2.Range("A4").Value = "Apple"
$ws2.Range("B4").Value = "MacBook Pro"
$ws2.Range("C4").Value = "Space Gray"
$ws2.Range("D4").Value = 2023

$ws2.Range("A5").Value = "Lenovo"
$ws2.Range("B5").Value = "ThinkPad X1 Carbon"
$ws2.Range("C5").Value = "Black"
$ws2.Range("D5").Value = 2022

$ws2.Range("A6").Value = "Acer"
$ws2.Range("B6").Value = "Swift 3"
$ws2.Range("C6").Value = "Gray"
$ws2.Range("D6").Value = 2020

$ws2.Columns.Item(2).ColumnWidth = 18.26953125
$ws2.Columns.Item(3).ColumnWidth = 11.453125

$ws2.Range("E3").Select()

$ws1.Select()
$ws1.Range("F11:F12").Select()
$ws1.Range("F12").Activate()
